$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value2 = -4.346463241969856
$ws.Range("C2").Value2 = 0.001136469920403883
$ws.Range("D2").Value2 = 0.01000230000030109
$ws.Range("E2").Value2 = 0.4320365356713142
$ws.Range("L2").Value2 = 0.007546198356849778
$ws.Range("M2").Value2 = 0.1141431108761355
$ws.Range("N2").Value2 = 0.00005472979131896184
$ws.Range("O2").Value2 = 0.1141431108761355
$ws.Range("P2").Value2 = 1.305964121354403
$ws.Range("Q2").Value2 = 0.3722524839430216
$ws.Range("T2").Value2 = 0.2546807723500223
$ws.Range("U2").Value2 = 0.3650885884167779
$ws.Range("P3").Value2 = 15.28067713828226
$ws.Range("Q3").Value2 = 0.7630366084239169
$ws.Range("T3").Value2 = 0.9018581552247483
$ws.Range("U3").Value2 = 2.847271312111027
$ws.Range("D4").Value2 = 0.009858725955514484
$ws.Range("E4").Value2 = 0.452704544934202
$ws.Range("F4").Value2 = 0.007177060046642438
$ws.Range("G4").Value2 = 0.006099549357745674
$ws.Range("J4").Value2 = 2.065385710343613
$ws.Range("K4").Value2 = 0.4821551753572025
$ws.Range("P4").Value2 = 1.309590237555901
$ws.Range("Q4").Value2 = 0.2309107350701574
$ws.Range("R4").Value2 = 0.1951177236046016
$ws.Range("S4").Value2 = 0.4863252868335161
$ws.Range("T4").Value2 = 0.254451966892874
$ws.Range("U4").Value2 = 0.2769718537359103
$ws.Range("F5").Value2 = 0.007173967396683234
$ws.Range("G5").Value2 = 0.005686529708348604
$ws.Range("J5").Value2 = 5.719299224597668
$ws.Range("K5").Value2 = 2.244714942757736
$ws.Range("P5").Value2 = 3.607033867293345
$ws.Range("Q5").Value2 = 0.4398563942857729
$ws.Range("T5").Value2 = 0.02107338269331508
$ws.Range("U5").Value2 = 0.2154674203966001
$ws.Range("J6").Value2 = 2.058691461260213
$ws.Range("K6").Value2 = 0.4290959512891188
$ws.Range("P6").Value2 = 1.306144503106202
$ws.Range("Q6").Value2 = 0.3270279760509027
$ws.Range("T6").Value2 = 0.2542348737881068
$ws.Range("U6").Value2 = 0.2907207943990641
$ws.Range("J7").Value2 = 18.71853147365436
$ws.Range("K7").Value2 = 3.579091055796329
$ws.Range("P7").Value2 = 11.51100166555282
$ws.Range("Q7").Value2 = 0.379460808073765
$ws.Range("T7").Value2 = 1.133871700599248
$ws.Range("U7").Value2 = 1.401243056938023
$ws.Range("H8").Value2 = 2.699364521708428
$ws.Range("I8").Value2 = 0.007298460639385412
$ws.Range("P8").Value2 = 1.309935559217266
$ws.Range("Q8").Value2 = 0.302136297629607
$ws.Range("T8").Value2 = 0.2555660680303823
$ws.Range("U8").Value2 = 0.2552495209300517
$ws.Range("J9").Value2 = 26.44082220731349
$ws.Range("K9").Value2 = 3.109065176893555
$ws.Range("P9").Value2 = 15.78362657429639
$ws.Range("Q9").Value2 = 0.4125288166592271
$ws.Range("T9").Value2 = 0.1766721985308917
$ws.Range("U9").Value2 = 0.351597563992787
$ws.Range("H10").Value2 = 2.699222237276266
$ws.Range("I10").Value2 = 0.006899427687144034
$ws.Range("P10").Value2 = 1.305375885283782
$ws.Range("Q10").Value2 = 0.2872659353666364
$ws.Range("T10").Value2 = 0.2551266661437217
$ws.Range("U10").Value2 = 0.2919022498101783
$ws.Range("D11").Value2 = 0.009617090543235918
$ws.Range("E11").Value2 = 2.503589336580748
$ws.Range("H11").Value2 = 31.71856395076951
$ws.Range("I11").Value2 = 0.02748518079900007
$ws.Range("J11").Value2 = 24.85144455077353
$ws.Range("K11").Value2 = 2.512834449749025
$ws.Range("P11").Value2 = 15.02000665612014
$ws.Range("Q11").Value2 = 0.2608301500611638
$ws.Range("T11").Value2 = 0.007330485958231978
$ws.Range("U11").Value2 = 0.06275040922461911
$ws.Range("H12").Value2 = 2.699528621760853
$ws.Range("I12").Value2 = 0.005804521755796589
$ws.Range("P12").Value2 = 1.3076281020285
$ws.Range("Q12").Value2 = 0.215815776003293
$ws.Range("T12").Value2 = 0.2542116114137046
$ws.Range("U12").Value2 = 0.2210895279536588
$ws.Range("D13").Value2 = 0.008631371408935825
$ws.Range("E13").Value2 = 3.135330266206472
$ws.Range("F13").Value2 = 0.007174134149524174
$ws.Range("G13").Value2 = 0.003600068728041768
$ws.Range("H13").Value2 = 19.64684962058278
$ws.Range("I13").Value2 = 0.02400726085605208
$ws.Range("P13").Value2 = 9.313938542111591
$ws.Range("Q13").Value2 = 0.3585395541305567
$ws.Range("T13").Value2 = 0.2083871059407842
$ws.Range("U13").Value2 = 0.3911703803145124
$ws.Range("D14").Value2 = 0.009885461325426443
$ws.Range("E14").Value2 = 0.4626693404894814
$ws.Range("F14").Value2 = 0.0071757431980022
$ws.Range("G14").Value2 = 0.004814947201412341
$ws.Range("P14").Value2 = 1.304506851212454
$ws.Range("Q14").Value2 = 0.2363366977725388
$ws.Range("T14").Value2 = 0.2532549146591672
$ws.Range("U14").Value2 = 0.3810547960272023
$ws.Range("F15").Value2 = 0.007172717586218937
$ws.Range("G15").Value2 = 0.003021117408196237
$ws.Range("J15").Value2 = 20.55111579011857
$ws.Range("K15").Value2 = 2.358604827910589
$ws.Range("P15").Value2 = 12.44650679596932
$ws.Range("Q15").Value2 = 0.2349775634081312
$ws.Range("T15").Value2 = 0.01050847952531234
$ws.Range("U15").Value2 = 0.05844981121142693
$ws.Range("D16").Value2 = 0.009899606839474754
$ws.Range("E16").Value2 = 0.4168767132109167
$ws.Range("H16").Value2 = 2.699349154270508
$ws.Range("I16").Value2 = 0.006743728088740613
$ws.Range("P16").Value2 = 1.306724378472846
$ws.Range("Q16").Value2 = 0.2950921642076711
$ws.Range("T16").Value2 = 0.2543400472003792
$ws.Range("U16").Value2 = 0.2769718537359103
